# DI_Timer_App/TimerDisplayTable.xlsx
# Commit: "Added Wrathborne invasion."
# Append a new event row (A12) with the shared string "Wrathborne Invasion 12PM"
# and leave the final selection where Excel would land after the edit (F14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of event data, directly below the existing "Ancient Arena 930PM" row.
$ws.Range("A12").Value = "Wrathborne Invasion 12PM"

# Matches the post-edit cursor position recorded in the saved workbook.
$ws.Range("F14").Select()
